$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Add-TestRow {
    param(
        [string]$CaseText,
        [string[]]$BehaviorParts,
        [string]$SuccessText
    )

    $row = $t.Rows.Add()

    $caseCell = $row.Cells.Item(1)
    $caseCell.Range.Text = $CaseText

    $behaviorCell = $row.Cells.Item(2)
    $behaviorCell.Range.Text = $BehaviorParts[0]
    for ($i = 1; $i -lt $BehaviorParts.Count; $i++) {
        $endPos = $behaviorCell.Range.End
        $insertionPoint = $d.Range($endPos - 1, $endPos - 1)
        $insertionPoint.InsertAfter($BehaviorParts[$i])
    }

    $successCell = $row.Cells.Item(3)
    $successCell.Range.Text = $SuccessText
}

Add-TestRow "Accessing Settings Menu" `
    @("Settings menu shows when the settings button is pressed and scales properly to different iPhone screen sizes", ".") `
    "Y"

Add-TestRow "Changing Schedule" `
    @("Changes to schedule show after ", "using the settings menu to make the changes.") `
    "Y"

Add-TestRow "Log out Button" `
    @("Pressing the log out button within the settings menu returns the user to the login page") `
    "Y"
